$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing existing row 10 down to row 11
$ws.Rows.Item(10).Insert()

# New row 10 content
$ws.Range("A10").Value = "Ruilin"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "无"
$ws.Range("D10").Value = "QSN"
$ws.Range("E10").Value = "RES"
$ws.Range("F10").Value = "a5228610-fe6d-4383-b598-a7c34c3b8714"
$ws.Range("G10").Value = "HyRnez-RW_annotated.xlsx"
$ws.Range("H10").Value = "Why is this result not compared to in Table 1?"

# Row 11 now has the original data (moved down); update D11/E11/F11/G11/H11 to new values
$ws.Range("D11").Value = "DIS"
$ws.Range("E11").Value = "MET"
$ws.Range("F11").Value = "d3fb2dcb-ee08-4432-9f4b-c252dbb3433f"
$ws.Range("G11").Value = "SJ3dBGZ0Z_annotated.xlsx"
$ws.Range("H11").Value = "We evaluate our method on NLP task for two reasons: 1) they are particularly well-suited for evaluating our method (naturally large output spaces) 2) we did not dispose of the computational resources to tackle tasks from other domains such as vision (e.g. Flickr100M) which requires hundreds of GPUs for weeks."
